# New project V3.1.0 - Nuevo metodo para ejecutar un caso con varios datos
# Change the browser value in row 4 (Test Case Name3) from "firefox" to "chrome"
# and update the active selection to B4.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value = "chrome"

$ws.Range("B4").Select()
